$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# "Implement attacks" row: Effort increased from 8 to 12 (all attacks finished,
# excluding Combo and Soulbar attacks). The "Remaining" column (E14) is a
# shared formula (=C14-D14) and recalculates automatically (15-12=3).
$ws.Range("D14").Value = 12

# Update the view: scroll so row 3 is the top-visible row, and move the
# active selection to D15.
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D15").Select()
